$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: date 44661 -> 44662, new time label, hours 2 -> 5.5, add activity note ---
$ws.Range("A12").Value = 44662
$ws.Range("B12").Value = "9.30 - 15.00"
$ws.Range("D12").Value = 5.5
$ws.Range("E12").Value = "Enemy FSM, sound"

# --- Row 13: date 44662 -> 44663, time label -> existing "9.00 - 14.00", hours 2.5 -> 5, add note ---
$ws.Range("A13").Value = 44663
$ws.Range("B13").Value = "9.00 - 14.00"
$ws.Range("D13").Value = 5
$ws.Range("E13").Value = "Level design"

# --- Row 14: date 44663 -> 44664, add note (time/hours unchanged) ---
$ws.Range("A14").Value = 44664
$ws.Range("E14").Value = "Level design"

# --- Row 15: date 44664 -> 44665, new time label, hours 5 -> 2.5, new activity (was "9.00 - 14.00") ---
$ws.Range("A15").Value = 44665
$ws.Range("B15").Value = "10.00  - 12.30"
$ws.Range("D15").Value = 2.5
$ws.Range("E15").Value = "Level design, background music player"

# --- Row 16: date 44665 -> 44699, hours 5 -> 3; note set before the time label to match string order ---
$ws.Range("A16").Value = 44699
$ws.Range("E16").Value = "Continuing for redo. Refining enemy behaviour"
$ws.Range("B16").Value = "9.00 - 12.00"
$ws.Range("D16").Value = 3

# --- Row 17: date 44666 -> 44700, new time label, hours 5 -> 2, new note ---
$ws.Range("A17").Value = 44700
$ws.Range("B17").Value = "14.30 - 16.30"
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = " Working on gunplay"

# --- Rows 18-20: clear all data but keep the date-formatted A cell ---
$ws.Range("A18:E18").ClearContents()
$ws.Range("A19:E19").ClearContents()
$ws.Range("A20:E20").ClearContents()

# --- Row 21: fully removed from the sheet ---
$ws.Range("A21:E21").EntireRow.Delete()

# --- Total row: was row 30 (=SUM(D3:D29)), now row 29 (=SUM(D3:D28)) ---
$ws.Range("D30").ClearContents()
$ws.Range("D29").Formula = "=SUM(D3:D28)"

# --- Header row 2 (second table): Hours/Activity shift right one column (Q2 -> blank, R2/S2 filled) ---
$ws.Range("Q2").ClearContents()
$ws.Range("R2").Value = "Hours"
$ws.Range("S2").Value = "Activity"

# --- New formatted (date-style) blank cell introduced in row 3 ---
$ws.Range("O3").NumberFormat = $ws.Range("A3").NumberFormat

# --- Selection shown in the saved file ---
$ws.Range("K26").Select()
